$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (bus 632): update B2 only
$ws.Range("B2").Value = "1.0000 at 0.00"

# Row 3 (was bus 645, now bus 633)
$ws.Range("A3").Value = 633
$ws.Range("B3").Value = "0.9969 at -0.07"
$ws.Range("C3").Value = "0.9980 at -120.05"
$ws.Range("D3").Value = "0.9973 at 119.99"

# Row 4 (was bus 646, now bus 634)
$ws.Range("A4").Value = 634
$ws.Range("B4").Value = "0.9724 at -0.77"
$ws.Range("C4").Value = "0.9788 at -120.55"
$ws.Range("D4").Value = "0.9782 at 119.50"

# Row 5 (new, bus 671)
$ws.Range("A5").Value = 671
$ws.Range("B5").Value = "1.0011 at 0.01"
$ws.Range("C5").Value = "0.9980 at -120.01"
$ws.Range("D5").Value = "0.9973 at 119.81"
